# Add a new 'Correction ' column (N) to the Card10 sheet, right after the
# existing 'Event' column (M), and normalize the 'Event' header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# Fix header text in M1: "Event " (trailing space) -> "Event"
$ws.Range("M1").Value = "Event"

# New header in N1
$ws.Range("N1").Value = "Correction "

# Copy the header style (bold/border/centered) from M1 onto N1
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill data rows 2-13.
# Column M (currently blank inlineStr cells) gets "nan" for every row.
# Column N (new) stays blank for every row.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 13).Value = "nan"   # column M = 13
    $ws.Cells.Item($row, 14).Value = ""      # column N = 14
}
